# Rename the first three worksheets to the new SCORECrossSect_* names
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("STUDY_AllVariables").Name = "SCORECrossSect_AllVariables"
$wb.Worksheets.Item("VariableMapping").Name = "SCORECrossSect_VariableMapping"
$wb.Worksheets.Item("ValueMapping").Name = "SCORECrossSect_ValueMapping"

# Move the active tab / selection from VariableMapping over to Notes
$wb.Worksheets.Item("Notes").Activate()
